$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Appends one new record row (row 6) to the PRINCIPAL sheet's data table.
# Leading apostrophes force literal-text storage (including for the purely
# numeric FRU code and the three intentionally blank sub-category columns),
# matching the existing text-typed columns already used by rows 2-5.
$row6 = $ws.Range("A6:M6")

$ws.Range("A6").Value = "'DF"
$ws.Range("B6").Value = "'1234567"
$ws.Range("C6").Value = "'"
$ws.Range("D6").Value = "'"
$ws.Range("E6").Value = "'"
$ws.Range("F6").Value = "'TESTE"
$ws.Range("G6").Value = "'DS8K"
$ws.Range("H6").Value = "'SICOOB - (78KKT90 14/11/25_24/7) - DF"
$ws.Range("I6").Value = "'14/11/25"
$ws.Range("J6").Value = "'24/7"
$ws.Range("K6").Value = "'14/11/25"
$ws.Range("L6").Value = "'DENTRO"
$ws.Range("M6").Value = "'"

# The quote-prefix entry above stamps an implicit "quotePrefix" style on the
# cells; put the row back on the workbook's default (unstyled) formatting so
# it matches the plain data rows above it (rows 2-5 carry no explicit style).
$row6.Style = "Normal"
